# [#1149] isoformat context dates in spreadsheet style error tables
#
# Insert a new "acceptable" example row into the grants sheet (with a real
# date value formatted as an ISO YYYY-MM-DD date) and tidy up the sheet
# selections/active-tab bookkeeping to match.

$wb = $excel.ActiveWorkbook

$grants = $wb.Worksheets.Item("grants")
$extra  = $wb.Worksheets.Item("extra sheet")

# 1. Make room for a new "acceptable" example row just under the header.
$grants.Rows.Item(2).Insert()

# 2. Populate the new row with a fully valid example grant.
$grants.Range("A2").Value = "360G-sampletrust-105177/Z/13/Z"
$grants.Range("B2").Value = "Acceptable title"
$grants.Range("C2").Value = "Acceptable description"
$grants.Range("D2").Value = "GBP"
$grants.Range("E2").Value = 1000

# Award Date: a real date value, isoformatted via a custom number format.
$grants.Range("F2").Value = 43617
$grants.Range("F2").NumberFormat = "YYYY\-MM\-DD"

# 3. Widen the new Award Date column slightly so the ISO date fits.
$grants.Columns.Item(6).ColumnWidth = 10.17

# 4. Refresh the selections: "extra sheet" keeps A10 selected (no longer the
#    active tab), "grants" becomes the active tab with F3 selected.
$extra.Range("A10").Select()
$grants.Activate()
$grants.Range("F3").Select()
